# Reduced cost and mass FELA
# - Fix typo "Phased Array (Deploable)" -> "Phased Array (Deployable)" in A3
# - Add three new rows (10-12) for "Multi Bandwidth Dish Transceiver" variants

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in A3 (Phased Array Deployable)
$ws.Range("A3").Value = "Phased Array (Deployable)"

# New row 10: Multi Bandwidth Dish Transceiver (large)
$ws.Range("A10").Value = "Multi Bandwidth Dish Transceiver (large)"
$ws.Range("C10").Value = 20
$ws.Range("F10").Value = "yes"
$ws.Range("G10").Value = 24
$ws.Range("H10").Value = 21000

# New row 11: Multi Bandwidth Dish Transceiver (medium)
$ws.Range("A11").Value = "Multi Bandwidth Dish Transceiver (medium)"
$ws.Range("C11").Value = 10
$ws.Range("F11").Value = "yes"
$ws.Range("G11").Value = 8
$ws.Range("H11").Value = 7000

# New row 12: Multi Bandwidth Dish Transceiver (shielded)
$ws.Range("A12").Value = "Multi Bandwidth Dish Transceiver (shielded)"
$ws.Range("C12").Value = 5
$ws.Range("F12").Value = "yes"
$ws.Range("G12").Value = 4
$ws.Range("H12").Value = 3500

# Update selection to match the authored workbook state
$ws.Range("L6").Select() | Out-Null
